# Apply reporting log updates to the trapping report worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that already had a species recorded (MARTEN) just need re-casing.
$ws.Range("F2").Value = "Marten"
$ws.Range("F5").Value = "Marten"
$ws.Range("F8").Value = "Marten"

# Rows that previously recorded "No harvest" now record a harvest of an
# unknown/"Na" species with one individual of unknown sex.
$noHarvestRows = @(3, 4, 6, 7, 9, 10)
foreach ($row in $noHarvestRows) {
    $ws.Range("B$row").Value = "Yes"
    $ws.Range("F$row").Value = "Na"
    $ws.Range("J$row").Value = 1
}
